$wb = $excel.ActiveWorkbook

# --- Sheet3 (Waiting-on log) ------------------------------------------------
# Rabbi Bordon's item has been resolved; replace it with the new entry for
# Rabbi Klein (waiting on Him), refresh the "As at" date, and drop the old
# duplicate Rabbi Klein / Him row that's no longer needed underneath.
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("I1").Value = "As at 25/02/19"
$ws3.Range("I3").Value = "Rabbi Klein"
$ws3.Range("J3").Value = "Him"
$ws3.Range("I4").ClearContents()
$ws3.Range("J4").ClearContents()

# Move the active selection to B27 (next free row).
$ws3.Range("B27").Select()

# --- Sheet1 ("STATUS") -----------------------------------------------------
# Mark H10 (Second review) and F15 (Reviewed?) as complete (TRUE).
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("H10").Value = $true
$ws1.Range("F15").Value = $true

# Move the active selection to D16 (reflects where the editor was working).
# Selected last so Sheet1 stays the active/tab-selected sheet, matching the
# workbook's original tab state.
$ws1.Range("D16").Select()
